$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 37.84230566666667
$ws.Range("H2").Value = 113.526917
$ws.Range("I2").Value = 0.5048163536019187
$ws.Range("J2").Value = 0.5048163536019187
$ws.Range("M2").Value = 0.02308233333333333
$ws.Range("N2").Value = 0.069247
$ws.Range("O2").Value = 0.003332618012635044
$ws.Range("P2").Value = 0.003332618012635044
$ws.Range("Q2").Value = 0.8734887134998889
$ws.Range("R2").Value = 7.861398421499
$ws.Range("S2").Value = 0.001682360073086496
$ws.Range("T2").Value = 0.001682360073086496
$ws.Range("G3").Value = 37.84230566666667
$ws.Range("H3").Value = 113.526917
$ws.Range("I3").Value = 0.5048163536019187
$ws.Range("J3").Value = 0.5048163536019187
$ws.Range("O3").Value = 0.03125495853682635
$ws.Range("P3").Value = 0.03125495853682635
$ws.Range("Q3").Value = 8.192014032006778
$ws.Range("R3").Value = 73.72812628806101
$ws.Range("S3").Value = 0.01577801420053984
$ws.Range("T3").Value = 0.01577801420053984
$ws.Range("G4").Value = 37.84230566666667
$ws.Range("H4").Value = 113.526917
$ws.Range("I4").Value = 0.5048163536019187
$ws.Range("J4").Value = 0.5048163536019187
$ws.Range("M4").Value = 0.055989
$ws.Range("N4").Value = 0.167967
$ws.Range("O4").Value = 0.008083669324711114
$ws.Range("P4").Value = 0.008083669324711114
$ws.Range("Q4").Value = 2.118752851971
$ws.Range("R4").Value = 19.068775667739
$ws.Range("S4").Value = 0.004080768472224348
$ws.Range("T4").Value = 0.004080768472224348
$ws.Range("G5").Value = 37.84230566666667
$ws.Range("H5").Value = 113.526917
$ws.Range("I5").Value = 0.5048163536019187
$ws.Range("J5").Value = 0.5048163536019187
$ws.Range("M5").Value = 6.630637333333333
$ws.Range("N5").Value = 19.891912
$ws.Range("O5").Value = 0.9573287541258274
$ws.Range("P5").Value = 0.9573287541258275
$ws.Range("Q5").Value = 250.9186047328115
$ws.Range("R5").Value = 2258.267442595304
$ws.Range("S5").Value = 0.4832752108560679
$ws.Range("T5").Value = 0.483275210856068
$ws.Range("G6").Value = 15.88630666666666
$ws.Range("H6").Value = 47.65891999999999
$ws.Range("I6").Value = 0.2119233292577262
$ws.Range("J6").Value = 0.2119233292577262
$ws.Range("M6").Value = 0.02308233333333333
$ws.Range("N6").Value = 0.069247
$ws.Range("O6").Value = 0.003332618012635044
$ws.Range("P6").Value = 0.003332618012635044
$ws.Range("Q6").Value = 0.3666930259155555
$ws.Range("R6").Value = 3.30023723324
$ws.Range("S6").Value = 0.0007062595043818854
$ws.Range("T6").Value = 0.0007062595043818856
$ws.Range("G7").Value = 15.88630666666666
$ws.Range("H7").Value = 47.65891999999999
$ws.Range("I7").Value = 0.2119233292577262
$ws.Range("J7").Value = 0.2119233292577262
$ws.Range("O7").Value = 0.03125495853682635
$ws.Range("P7").Value = 0.03125495853682635
$ws.Range("Q7").Value = 3.439030599151111
$ws.Range("R7").Value = 30.95127539236
$ws.Range("S7").Value = 0.00662365486893643
$ws.Range("T7").Value = 0.006623654868936431
$ws.Range("G8").Value = 15.88630666666666
$ws.Range("H8").Value = 47.65891999999999
$ws.Range("I8").Value = 0.2119233292577262
$ws.Range("J8").Value = 0.2119233292577262
$ws.Range("M8").Value = 0.055989
$ws.Range("N8").Value = 0.167967
$ws.Range("O8").Value = 0.008083669324711114
$ws.Range("P8").Value = 0.008083669324711114
$ws.Range("Q8").Value = 0.88945842396
$ws.Range("R8").Value = 8.00512581564
$ws.Range("S8").Value = 0.001713118115911334
$ws.Range("T8").Value = 0.001713118115911335
$ws.Range("G9").Value = 15.88630666666666
$ws.Range("H9").Value = 47.65891999999999
$ws.Range("I9").Value = 0.2119233292577262
$ws.Range("J9").Value = 0.2119233292577262
$ws.Range("M9").Value = 6.630637333333333
$ws.Range("N9").Value = 19.891912
$ws.Range("O9").Value = 0.9573287541258274
$ws.Range("P9").Value = 0.9573287541258275
$ws.Range("Q9").Value = 105.3363380727822
$ws.Range("R9").Value = 948.0270426550398
$ws.Range("S9").Value = 0.2028802967684965
$ws.Range("T9").Value = 0.2028802967684966
$ws.Range("G10").Value = 18.76675533333333
$ws.Range("H10").Value = 56.300266
$ws.Range("I10").Value = 0.2503485141672444
$ws.Range("J10").Value = 0.2503485141672445
$ws.Range("M10").Value = 0.02308233333333333
$ws.Range("N10").Value = 0.069247
$ws.Range("O10").Value = 0.003332618012635044
$ws.Range("P10").Value = 0.003332618012635044
$ws.Range("Q10").Value = 0.4331805021891111
$ws.Range("R10").Value = 3.898624519702
$ws.Range("S10").Value = 0.0008343159677501781
$ws.Range("T10").Value = 0.0008343159677501785
$ws.Range("G11").Value = 18.76675533333333
$ws.Range("H11").Value = 56.300266
$ws.Range("I11").Value = 0.2503485141672444
$ws.Range("J11").Value = 0.2503485141672445
$ws.Range("O11").Value = 0.03125495853682635
$ws.Range("P11").Value = 0.03125495853682635
$ws.Range("Q11").Value = 4.062583405464222
$ws.Range("R11").Value = 36.563250649178
$ws.Range("S11").Value = 0.007824632430053309
$ws.Range("T11").Value = 0.00782463243005331
$ws.Range("G12").Value = 18.76675533333333
$ws.Range("H12").Value = 56.300266
$ws.Range("I12").Value = 0.2503485141672444
$ws.Range("J12").Value = 0.2503485141672445
$ws.Range("M12").Value = 0.055989
$ws.Range("N12").Value = 0.167967
$ws.Range("O12").Value = 0.008083669324711114
$ws.Range("P12").Value = 0.008083669324711114
$ws.Range("Q12").Value = 1.050731864358
$ws.Range("R12").Value = 9.456586779222
$ws.Range("S12").Value = 0.002023734604460759
$ws.Range("T12").Value = 0.00202373460446076
$ws.Range("G13").Value = 18.76675533333333
$ws.Range("H13").Value = 56.300266
$ws.Range("I13").Value = 0.2503485141672444
$ws.Range("J13").Value = 0.2503485141672445
$ws.Range("M13").Value = 6.630637333333333
$ws.Range("N13").Value = 19.891912
$ws.Range("O13").Value = 0.9573287541258274
$ws.Range("P13").Value = 0.9573287541258275
$ws.Range("Q13").Value = 124.4355485387324
$ws.Range("R13").Value = 1119.919936848592
$ws.Range("S13").Value = 0.2396658311649801
$ws.Range("T13").Value = 0.2396658311649802
$ws.Range("G14").Value = 2.467151666666667
$ws.Range("H14").Value = 7.401455
$ws.Range("I14").Value = 0.03291180297311068
$ws.Range("J14").Value = 0.03291180297311068
$ws.Range("M14").Value = 0.02308233333333333
$ws.Range("N14").Value = 0.069247
$ws.Range("O14").Value = 0.003332618012635044
$ws.Range("P14").Value = 0.003332618012635044
$ws.Range("Q14").Value = 0.05694761715388889
$ws.Range("R14").Value = 0.512528554385
$ws.Range("S14").Value = 0.0001096824674164842
$ws.Range("T14").Value = 0.0001096824674164843
$ws.Range("G15").Value = 2.467151666666667
$ws.Range("H15").Value = 7.401455
$ws.Range("I15").Value = 0.03291180297311068
$ws.Range("J15").Value = 0.03291180297311068
$ws.Range("O15").Value = 0.03125495853682635
$ws.Range("P15").Value = 0.03125495853682635
$ws.Range("Q15").Value = 0.5340832361127777
$ws.Range("R15").Value = 4.806749125015
$ws.Range("S15").Value = 0.001028657037296773
$ws.Range("T15").Value = 0.001028657037296773
$ws.Range("G16").Value = 2.467151666666667
$ws.Range("H16").Value = 7.401455
$ws.Range("I16").Value = 0.03291180297311068
$ws.Range("J16").Value = 0.03291180297311068
$ws.Range("M16").Value = 0.055989
$ws.Range("N16").Value = 0.167967
$ws.Range("O16").Value = 0.008083669324711114
$ws.Range("P16").Value = 0.008083669324711114
$ws.Range("Q16").Value = 0.138133354665
$ws.Range("R16").Value = 1.243200191985
$ws.Range("S16").Value = 0.0002660481321146708
$ws.Range("T16").Value = 0.0002660481321146708
$ws.Range("G17").Value = 2.467151666666667
$ws.Range("H17").Value = 7.401455
$ws.Range("I17").Value = 0.03291180297311068
$ws.Range("J17").Value = 0.03291180297311068
$ws.Range("M17").Value = 6.630637333333333
$ws.Range("N17").Value = 19.891912
$ws.Range("O17").Value = 0.9573287541258274
$ws.Range("P17").Value = 0.9573287541258275
$ws.Range("Q17").Value = 16.35878794799555
$ws.Range("R17").Value = 147.22909153196
$ws.Range("S17").Value = 0.03150741533628275
$ws.Range("T17").Value = 0.03150741533628276
